# Applies the "accounting for offset of additionality and impact" fix:
#   1. Row 8 (COLLE RENTAL & SALES ELECTRIFICATION) had its
#      additionality_and_impact / environmental_aspects / procurement values
#      shifted by one column (Q/R/S held "Procurement" / "Link to source" /
#      "Summary sheet" placeholders instead of real content). Fix the offset:
#      Q8 -> "N/A", R8 -> real environmental text, S8 -> real procurement text.
#   2. Append a brand-new row 9 for project KAYRROS ARTIFICIAL INTELLIGENCE
#      (EGF VD), with its own source hyperlink in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix the additionality/environmental/procurement offset on row 8
# ---------------------------------------------------------------------------
$ws.Range("Q8").Value = "N/A"

$ws.Range("R8").Value = @'
The project concerns capital expenditures on new rental equipment; the rental activities will be carried out in existing facilities already authorised for the same purpose and would not require an environmental impact assessment (EIA) under the Directive 2014/52/EU.
'@

$ws.Range("S8").Value = @'
The Promoter has been assessed by EIB as being a private company not operating in the utilities sector and not having a status of a contracting entity, and is thus not subject to EU rules on public procurement.
'@

# ---------------------------------------------------------------------------
# 2. Append new row 9 - KAYRROS ARTIFICIAL INTELLIGENCE (EGF VD)
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "https://www.eib.org/en/projects/all/20200483"
$ws.Hyperlinks.Add($ws.Range("A9"), "https://www.eib.org/en/projects/all/20200483")
$ws.Range("A9").Style = "Hyperlink"

$ws.Range("B9").Value = "28 October 2020"
$ws.Range("C9").Value = "Signed"
$ws.Range("D9").Value = "28/02/2022"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "20200483"
$ws.Range("E9").Style = "Normal"

$ws.Range("F9").Value = "KAYRROS ARTIFICIAL INTELLIGENCE (EGF VD)"
$ws.Range("G9").Value = "KAYRROS SAS"
$ws.Range("H9").Value = "EUR"

$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "34"
$ws.Range("I9").Style = "Normal"

$ws.Range("J9").Value = "million"
$ws.Range("K9").Value = "EUR"

$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "15"
$ws.Range("L9").Style = "Normal"

$ws.Range("M9").Value = "million"
$ws.Range("N9").Value = "France"

$ws.Range("O9").Value = @'
Kayrros is an asset observation platform powered by modern technologies such as data analytics, artificial intelligence (AI) and computer vision, that uses several data sources, including satellite images, Internet of Things (IoT) devices, sensors, mobile devices, geolocation, social media, web scraping, etc. to provide intelligence to private and public decision-makers.
'@

$ws.Range("P9").Value = @'
The project will finance investments for research, development and innovation (RDI) activities related to the development of data analytics and artificial intelligence, as well as on IT integration, to consolidate its first-mover advantage in the asset-observation data analytics market.
'@

$ws.Range("Q9").Value = @'
The promoter is an innovative start-up with well-developed technology capabilities, an organisational structure in support of its ambitious growth strategy, and a strong R&D team. The project will help the promoter pursue its R&D investments, accelerate the development and the deployment of new use cases, and build scale. The project supports the policy objective of innovation and
digitalisation in the asset observation domain, which has long been suffering
from imperfect and opaque information. By supporting the company's investments
in R&D and its expansion into new applications, including environment and
climate, risk management and the monitoring of the global supply chain, the
project also addresses knowledge externalities and contributes to public goods.
The promoter's strategy depends on significant growth-related investments in
order to achieve its business plan. By addressing a market gap in available financing
options, EIB financing allows the company to finance its investment plan and
accelerate the deployment of its growth strategy. EIB contribution is positive because of a positive crowd-in effect and innovative structure adapted
to the risk profile of this Company given the early stage of the company.
'@

$ws.Range("R9").Value = @'
The project activities do not fall under Annexes I or II of the EU Directive 2014/52/EU amending the EIA Directive 2011/92/EU. The project will be carried out in existing facilities, already authorised, that will not change their scope due to the project. As such, the project activities are not subject to a mandatory environmental impact assessment (EIA).
'@

$ws.Range("S9").Value = @'
The Promoter is a private company not operating in the utilities sector and not having the status of a contracting entity, and is thus not subject to EU rules on public procurement. Under these conditions, the procurement procedures followed by the Promoter are suitable for the project.
'@

$ws.Range("T9").Value = @'
KAYRROS ARTIFICIAL INTELLIGENCE (EGF VD): https://www.eib.org/en/projects/pipelines/all/20200483 || France: Climate technology specialist Kayrros receives €15 million EIB loan: https://www.eib.org/en/press/all/2022-138-france-kayrros-specialise-dans-la-tech-climatique-beneficie-d-un-pret-de-15-millions-d-euros-de-la-bei || 
'@

$ws.Range("U9").Value = "€"

$ws.Range("V9").NumberFormat = "@"
$ws.Range("V9").Value = "15,000,000"
$ws.Range("V9").Style = "Normal"

$ws.Range("W9").Value = "France"
$ws.Range("X9").Value = "€"

$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = "15,000,000"
$ws.Range("Y9").Style = "Normal"

$ws.Range("Z9:AQ9").Value = "No Entry"

$ws.Range("AR9").Value = "Services"
$ws.Range("AS9").Value = "Information and communication"
$ws.Range("AT9").Value = "€"

$ws.Range("AU9").NumberFormat = "@"
$ws.Range("AU9").Value = "15,000,000"
$ws.Range("AU9").Style = "Normal"

$ws.Range("AV9:BC9").Value = "No Entry"

$ws.Range("BD9").Value = $true
